# Re-pull / push data: update column F (dSF) values for several rows so
# that they match the freshly re-pulled source data (they previously had
# been left equal to column E / dS0 by mistake).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -4
    5  = 0
    6  = -6
    8  = -6
    9  = -2
    14 = -10
    20 = -6
    21 = 2
    25 = -5
    28 = -6
    34 = -8
    37 = -6
    39 = 0
    44 = 2
    47 = -4
    52 = 0
    53 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
